$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 2846.1538
$ws.Range("I54").Value = 9000
$ws.Range("J54").Value = 1000
$ws.Range("K54").Value = 9000
$ws.Range("L54").Value = 1000
$ws.Range("M54").Value = -8514
$ws.Range("N54").Value = -1972
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2434.0667
$ws.Range("I2").Value = 1649.8
$ws.Range("K2").Value = 1649.8
$ws.Range("M2").Value = -1536.8
$ws.Range("H88").Value = 9555.857
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 9555.857
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H92").Value = 34025
$ws.Range("I92").Value = 20000
$ws.Range("J92").Value = 38700
$ws.Range("K92").Value = 20000
$ws.Range("L92").Value = 38700
$ws.Range("M92").Value = -17504
$ws.Range("N92").Value = -43692
$ws.Range("H97").Value = 1282.4166
$ws.Range("I97").Value = 722.2222
$ws.Range("K97").Value = 722.2222
$ws.Range("M97").Value = -226.2222
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H116").Value = 2434.0667
$ws.Range("I116").Value = 1649.8
$ws.Range("K116").Value = 1649.8
$ws.Range("M116").Value = 644.2
$ws.Range("H130").Value = 12999.333
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 12999.333
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 12999.333
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -23039.333
$ws.Range("H132").Value = 4700.778
$ws.Range("I132").Value = 4800.875
$ws.Range("J132").Value = 3900
$ws.Range("K132").Value = 14402.625
$ws.Range("L132").Value = 11700
$ws.Range("M132").Value = -11872.625
$ws.Range("N132").Value = -16760

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2434.0667
$ws.Range("I3").Value = 1649.8
$ws.Range("K3").Value = 1649.8
$ws.Range("M3").Value = -1535.8
$ws.Range("H12").Value = 1126.25
$ws.Range("I12").Value = 102.5
$ws.Range("J12").Value = 2150
$ws.Range("K12").Value = 102.5
$ws.Range("L12").Value = 2150
$ws.Range("M12").Value = 65.5
$ws.Range("N12").Value = -2486
$ws.Range("H86").Value = 2003.5
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H89").Value = 2003.5
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("H99").Value = 5198.3335
$ws.Range("I99").Value = 5198.3335
$ws.Range("K99").Value = 5198.3335
$ws.Range("M99").Value = -3700.3335

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2714.0908
$ws.Range("I31").Value = 2169.4
$ws.Range("J31").Value = 3168
$ws.Range("K31").Value = 2169.4
$ws.Range("L31").Value = 3168
$ws.Range("M31").Value = -1874.4
$ws.Range("N31").Value = -3758
$ws.Range("H34").Value = 2714.0908
$ws.Range("I34").Value = 2169.4
$ws.Range("J34").Value = 3168
$ws.Range("K34").Value = 2169.4
$ws.Range("L34").Value = 3168
$ws.Range("M34").Value = -1967.4
$ws.Range("N34").Value = -3572
$ws.Range("H74").Value = 40000
$ws.Range("J74").Value = 40000
$ws.Range("L74").Value = 40000
$ws.Range("N74").Value = -41748
$ws.Range("H77").Value = 40000
$ws.Range("J77").Value = 40000
$ws.Range("L77").Value = 120000
$ws.Range("N77").Value = -128736
$ws.Range("H80").Value = 25200
$ws.Range("I80").Value = 25200
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 25200
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -24077
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 25200
$ws.Range("I83").Value = 25200
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 75600
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -69984
$ws.Range("N83").ClearContents()
$ws.Range("H86").Value = 3222
$ws.Range("I86").Value = 3222
$ws.Range("K86").Value = 3222
$ws.Range("M86").Value = -2099
$ws.Range("H89").Value = 3222
$ws.Range("I89").Value = 3222
$ws.Range("K89").Value = 16110
$ws.Range("M89").Value = -10494
$ws.Range("H134").Value = 3462.842
$ws.Range("I134").Value = 3458.7058
$ws.Range("K134").Value = 10376.1174
$ws.Range("M134").Value = -7841.117400000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 164.8
$ws.Range("I13").Value = 164.8
$ws.Range("K13").Value = 494.4
$ws.Range("M13").Value = -326.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 100777
$ws.Range("J133").Value = 100777
$ws.Range("L133").Value = 100777
$ws.Range("N133").Value = -110897

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H132").Value = 17276.375
$ws.Range("I132").Value = 18895.092
$ws.Range("J132").Value = 13715.2
$ws.Range("K132").Value = 56685.276
$ws.Range("L132").Value = 41145.60000000001
$ws.Range("M132").Value = -54155.276
$ws.Range("N132").Value = -46205.60000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 102500.5
$ws.Range("J46").Value = 102500.5
$ws.Range("L46").Value = 102500.5
$ws.Range("N46").Value = -102962.5
$ws.Range("H110").Value = 24999
$ws.Range("J110").Value = 24999
$ws.Range("L110").Value = 24999
$ws.Range("N110").Value = -33179
$ws.Range("H134").Value = 102500.5
$ws.Range("J134").Value = 102500.5
$ws.Range("L134").Value = 307501.5
$ws.Range("N134").Value = -312571.5
